$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2689.611
$ws.Range("I40").Value = 3883.3333
$ws.Range("J40").Value = 2092.75
$ws.Range("K40").Value = 3883.3333
$ws.Range("L40").Value = 2092.75
$ws.Range("M40").Value = -3708.3333
$ws.Range("N40").Value = -2442.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3486.0857
$ws.Range("I64").Value = 3321.8928
$ws.Range("K64").Value = 3321.8928
$ws.Range("M64").Value = -3073.8928

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3486.0857
$ws.Range("I67").Value = 3321.8928
$ws.Range("K67").Value = 3321.8928
$ws.Range("M67").Value = -2463.8928

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 3898.8857
$ws.Range("I80").Value = 631
$ws.Range("J80").Value = 7779.5
$ws.Range("K80").Value = 1893
$ws.Range("L80").Value = 23338.5
$ws.Range("M80").Value = -895
$ws.Range("N80").Value = -25334.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 3898.8857
$ws.Range("I83").Value = 631
$ws.Range("J83").Value = 7779.5
$ws.Range("K83").Value = 5679
$ws.Range("L83").Value = 70015.5
$ws.Range("M83").Value = -687
$ws.Range("N83").Value = -79999.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1582.3914
$ws.Range("I2").Value = 1594.9286
$ws.Range("J2").Value = 1562.8889
$ws.Range("K2").Value = 1594.9286
$ws.Range("L2").Value = 1562.8889
$ws.Range("M2").Value = -1481.9286
$ws.Range("N2").Value = -1788.8889

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 2303.4
$ws.Range("I14").Value = 379.25
$ws.Range("J14").Value = 10000
$ws.Range("K14").Value = 379.25
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = -204.25
$ws.Range("N14").Value = -10350

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2113.62
$ws.Range("I32").Value = 1985.0309
$ws.Range("J32").Value = 6271.3335
$ws.Range("K32").Value = 1985.0309
$ws.Range("L32").Value = 6271.3335
$ws.Range("M32").Value = -1698.0309
$ws.Range("N32").Value = -6845.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1900.9
$ws.Range("I63").Value = 1900.9
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 1900.9
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1214.9
$ws.Range("N63").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 1900.9
$ws.Range("I66").Value = 1900.9
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 9504.5
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -6072.5
$ws.Range("N66").Value = 0

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1582.3914
$ws.Range("I116").Value = 1594.9286
$ws.Range("J116").Value = 1562.8889
$ws.Range("K116").Value = 1594.9286
$ws.Range("L116").Value = 1562.8889
$ws.Range("M116").Value = 699.0714
$ws.Range("N116").Value = -6150.8889

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 7815626
$ws.Range("I122").Value = 3385.7693
$ws.Range("J122").Value = 41668668
$ws.Range("K122").Value = 10157.3079
$ws.Range("L122").Value = 125006004
$ws.Range("M122").Value = -7707.3079
$ws.Range("N122").Value = -125010904

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1582.3914
$ws.Range("I3").Value = 1594.9286
$ws.Range("J3").Value = 1562.8889
$ws.Range("K3").Value = 1594.9286
$ws.Range("L3").Value = 1562.8889
$ws.Range("M3").Value = -1480.9286
$ws.Range("N3").Value = -1790.8889

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 400
$ws.Range("I22").Value = 400
$ws.Range("K22").Value = 400
$ws.Range("M22").Value = -227

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2858.9836
$ws.Range("I134").Value = 2797.4385
$ws.Range("K134").Value = 8392.315500000001
$ws.Range("M134").Value = -5857.315500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1529.7273
$ws.Range("I94").Value = 1414.125
$ws.Range("J94").Value = 1838
$ws.Range("K94").Value = 1414.125
$ws.Range("L94").Value = 1838
$ws.Range("M94").Value = -963.125
$ws.Range("N94").Value = -2740

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 11649.5
$ws.Range("I122").Value = 3769.7856
$ws.Range("J122").Value = 39228.5
$ws.Range("K122").Value = 11309.3568
$ws.Range("L122").Value = 117685.5
$ws.Range("M122").Value = -8859.356800000001
$ws.Range("N122").Value = -122585.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3136.8333
$ws.Range("I134").Value = 2085.4285
$ws.Range("J134").Value = 4269.115
$ws.Range("K134").Value = 6256.2855
$ws.Range("L134").Value = 12807.345
$ws.Range("M134").Value = -3721.2855
$ws.Range("N134").Value = -17877.345

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 7684.933
$ws.Range("I87").Value = 3558.8
$ws.Range("J87").Value = 9748
$ws.Range("K87").Value = 10676.4
$ws.Range("L87").Value = 29244
$ws.Range("M87").Value = -9428.400000000001
$ws.Range("N87").Value = -31740

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 7684.933
$ws.Range("I90").Value = 3558.8
$ws.Range("J90").Value = 9748
$ws.Range("K90").Value = 32029.2
$ws.Range("L90").Value = 87732
$ws.Range("M90").Value = -25789.2
$ws.Range("N90").Value = -100212

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1456.0883
$ws.Range("J107").Value = 1859.96
$ws.Range("L107").Value = 5579.88
$ws.Range("N107").Value = -9419.880000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 10411.58
$ws.Range("I131").Value = 1067.381
$ws.Range("J131").Value = 12895.481
$ws.Range("K131").Value = 3202.143
$ws.Range("L131").Value = 38686.443
$ws.Range("M131").Value = 1837.857
$ws.Range("N131").Value = -48766.443

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 1000000
$ws.Range("I14").Value = 1000000
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 1000000
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -999832
$ws.Range("N14").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6650
$ws.Range("I122").Value = 7111.1113
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 21333.3339
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -18883.3339
$ws.Range("N122").Value = -12400

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 9033.789000000001
$ws.Range("J136").Value = 9033.789000000001
$ws.Range("L136").Value = 27101.367
$ws.Range("N136").Value = -32201.367

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3530.7727
$ws.Range("I100").Value = 1991.8
$ws.Range("J100").Value = 6828.5713
$ws.Range("K100").Value = 1991.8
$ws.Range("L100").Value = 6828.5713
$ws.Range("M100").Value = -1450.8
$ws.Range("N100").Value = -7910.5713

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2380
$ws.Range("I96").Value = 2306
$ws.Range("J96").Value = 2472.5
$ws.Range("K96").Value = 2306
$ws.Range("L96").Value = 2472.5
$ws.Range("M96").Value = -933
$ws.Range("N96").Value = -5218.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 945.2
$ws.Range("I107").Value = 351.11765
$ws.Range("J107").Value = 2207.625
$ws.Range("K107").Value = 1053.35295
$ws.Range("L107").Value = 6622.875
$ws.Range("M107").Value = 866.64705
$ws.Range("N107").Value = -10462.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3122.1333
$ws.Range("I122").Value = 1870.36
$ws.Range("K122").Value = 5611.08
$ws.Range("M122").Value = -3161.08
